$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Log Time" value in D2 was re-entered as a real Excel time value
# (serial fraction of a day) instead of plain text, and formatted with the
# built-in h:mm:ss time number format / a black font color.
$cell = $ws.Range("D2")
$cell.Value = 0.46059027777777778
$cell.NumberFormat = "h:mm:ss"
$cell.Font.Color = 0

# Touch row 3 (it now shows up - empty - below the data, as happens after
# Excel re-saves the sheet) while re-using D2's freshly created style so we
# don't introduce any extra/duplicate style entries.
$ws.Range("A3").Style = $cell.Style

# Row heights settle at Excel's normal 12pt-Calibri row height after resave.
$ws.Rows.Item(1).RowHeight = 15.5
$ws.Rows.Item(2).RowHeight = 15.5
$ws.Rows.Item(3).RowHeight = 15.5

# The active selection ends up sitting on D2.
$cell.Select()
